# Remove plate 1 from Suppl Table 1 (LAMP confusion matrix)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update raw count cells (formulas in column H / row 10 recalc automatically) ---
$ws.Range("E5").Value = 51
$ws.Range("E6").Value = 29
$ws.Range("E7").Value = 7
$ws.Range("G7").Value = 16
$ws.Range("G9").Value = 648

# --- Update statistic cells (plain stored values, not formulas) ---
$ws.Range("L5").Value = 0.929953380101468

$ws.Range("K6").Value = 0.852941176470588
$ws.Range("L6").Value = 0.698719435045478
$ws.Range("N6").Value = 0.935505585552175

$ws.Range("K7").Value = 0.304347826086957
$ws.Range("L7").Value = 0.15604024453214
$ws.Range("N7").Value = 0.508657562687592

$ws.Range("K12").Value = 0.996923076923077
$ws.Range("L12").Value = 0.988851303740987
$ws.Range("N12").Value = 0.999155790769604

# --- Row heights (rows 6,7,8,9,12 shrink from 17 to 15) ---
$ws.Rows.Item(6).RowHeight = 15
$ws.Rows.Item(7).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 15
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(12).RowHeight = 15

# --- View changes: zoom + selected cell ---
$excel.ActiveWindow.Zoom = 209
$ws.Range("D15").Select()
